$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (new iAuthor TC data)
$ws.Range("A2").Value = "AwKRq766"
$ws.Range("B2").Value = 23110739
$ws.Range("C2").Value = "bqetuen38"
$ws.Range("D2").Value = "M`$mY23!e"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "BPtatTwQ"
$ws.Range("G2").Value = "cjMy"
$ws.Range("H2").Value = "Candidate"

# Delete row 3 entirely (old second TC record removed)
$ws.Rows.Item(3).Delete()

# Keep the selection in sync with the new used range
$null = $ws.Range("A1:H2").Select()
